$wb = $excel.ActiveWorkbook

# Sheet ALC Row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 7747
$ws.Range("I86").Value = 6994
$ws.Range("K86").Value = 6994
$ws.Range("M86").Value = -5871

# Sheet ALC Row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 7747
$ws.Range("I89").Value = 6994
$ws.Range("K89").Value = 34970
$ws.Range("M89").Value = -29354

# Sheet ALC Row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 298.14285
$ws.Range("I96").Value = 298.14285
$ws.Range("K96").Value = 894.4285500000001
$ws.Range("M96").Value = 478.5714499999999

# Sheet ALC Row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3117.875
$ws.Range("J112").Value = 3117.875
$ws.Range("L112").Value = 9353.625
$ws.Range("N112").Value = -11569.625

# Sheet ALC Row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 679.0909
$ws.Range("I135").Value = 648
$ws.Range("J135").Value = 990
$ws.Range("K135").Value = 5832
$ws.Range("L135").Value = 8910
$ws.Range("M135").Value = -3297
$ws.Range("N135").Value = -13980

# Sheet ALC Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 35802.168
$ws.Range("I137").Value = 42073.36
$ws.Range("K137").Value = 126220.08
$ws.Range("M137").Value = -123670.08

# Sheet ALC Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3348.7568
$ws.Range("I138").Value = 4289.1333
$ws.Range("J138").Value = 2707.5908
$ws.Range("K138").Value = 12867.3999
$ws.Range("L138").Value = 8122.7724
$ws.Range("M138").Value = -7727.3999
$ws.Range("N138").Value = -18402.7724

# Sheet ARM Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23186.51
$ws.Range("I32").Value = 23924.791
$ws.Range("J32").Value = 16099
$ws.Range("K32").Value = 23924.791
$ws.Range("L32").Value = 16099
$ws.Range("M32").Value = -23637.791
$ws.Range("N32").Value = -16673

# Sheet ARM Row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1599.7
$ws.Range("I74").Value = 1678.8572
$ws.Range("J74").Value = 1415
$ws.Range("K74").Value = 1678.8572
$ws.Range("L74").Value = 1415
$ws.Range("M74").Value = -804.8571999999999
$ws.Range("N74").Value = -3163

# Sheet ARM Row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1599.7
$ws.Range("I77").Value = 1678.8572
$ws.Range("J77").Value = 1415
$ws.Range("K77").Value = 8394.286
$ws.Range("L77").Value = 7075
$ws.Range("M77").Value = -4026.286
$ws.Range("N77").Value = -15811

# Sheet ARM Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1909.2222
$ws.Range("I122").Value = 1909.2222
$ws.Range("K122").Value = 5727.6666
$ws.Range("M122").Value = -3277.6666

# Sheet ARM Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 30670.555
$ws.Range("I132").Value = 31204
$ws.Range("K132").Value = 93612
$ws.Range("M132").Value = -91082

# Sheet BSM Row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 558.8570999999999
$ws.Range("I80").Value = 836.6667
$ws.Range("J80").Value = 350.5
$ws.Range("K80").Value = 836.6667
$ws.Range("L80").Value = 350.5
$ws.Range("M80").Value = 161.3333
$ws.Range("N80").Value = -2346.5

# Sheet BSM Row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 558.8570999999999
$ws.Range("I83").Value = 836.6667
$ws.Range("J83").Value = 350.5
$ws.Range("K83").Value = 4183.3335
$ws.Range("L83").Value = 1752.5
$ws.Range("M83").Value = 808.6665000000003
$ws.Range("N83").Value = -11736.5

# Sheet BSM Row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4033.625
$ws.Range("I105").Value = 3824.1428
$ws.Range("K105").Value = 3824.1428
$ws.Range("M105").Value = -2077.1428

# Sheet BSM Row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3594.95
$ws.Range("I107").Value = 2854.4546
$ws.Range("J107").Value = 4500
$ws.Range("K107").Value = 2854.4546
$ws.Range("L107").Value = 4500
$ws.Range("M107").Value = -934.4546
$ws.Range("N107").Value = -8340

# Sheet CRP Row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 31561.914
$ws.Range("I58").Value = 60749.176
$ws.Range("J58").Value = 3996.1667
$ws.Range("K58").Value = 60749.176
$ws.Range("L58").Value = 3996.1667
$ws.Range("M58").Value = -60546.176
$ws.Range("N58").Value = -4402.1667

# Sheet CRP Row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6376.75
$ws.Range("I86").Value = 4749.5
$ws.Range("J86").Value = 8004
$ws.Range("K86").Value = 4749.5
$ws.Range("L86").Value = 8004
$ws.Range("M86").Value = -3626.5
$ws.Range("N86").Value = -10250

# Sheet CRP Row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 6376.75
$ws.Range("I89").Value = 4749.5
$ws.Range("J89").Value = 8004
$ws.Range("K89").Value = 23747.5
$ws.Range("L89").Value = 40020
$ws.Range("M89").Value = -18131.5
$ws.Range("N89").Value = -51252

# Sheet CRP Row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2012
$ws.Range("I99").Value = 2012
$ws.Range("K99").Value = 2012
$ws.Range("M99").Value = -514

# Sheet CRP Row 111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 29999
$ws.Range("J111").Value = 29999
$ws.Range("L111").Value = 29999
$ws.Range("N111").Value = -38179

# Sheet CRP Row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2012
$ws.Range("I126").Value = 2012
$ws.Range("K126").Value = 6036
$ws.Range("M126").Value = -3566

# Sheet CRP Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2592.0588
$ws.Range("I132").Value = 2050.0715
$ws.Range("J132").Value = 2971.45
$ws.Range("K132").Value = 6150.2145
$ws.Range("L132").Value = 8914.349999999999
$ws.Range("M132").Value = -3620.2145
$ws.Range("N132").Value = -13974.35

# Sheet CRP Row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 31561.914
$ws.Range("I136").Value = 60749.176
$ws.Range("J136").Value = 3996.1667
$ws.Range("K136").Value = 182247.528
$ws.Range("L136").Value = 11988.5001
$ws.Range("M136").Value = -179697.528
$ws.Range("N136").Value = -17088.5001

# Sheet CUL Row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 97.57143000000001
$ws.Range("J6").Value = 1
$ws.Range("L6").Value = 3
$ws.Range("N6").Value = -229

# Sheet CUL Row 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 190
$ws.Range("I50").Value = 190
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 570
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -89
$ws.Range("N50").ClearContents()

# Sheet CUL Row 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 190
$ws.Range("I53").Value = 190
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 570
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -89
$ws.Range("N53").ClearContents()

# Sheet CUL Row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2849.5
$ws.Range("J68").Value = 4999
$ws.Range("L68").Value = 14997
$ws.Range("N68").Value = -16619

# Sheet CUL Row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2849.5
$ws.Range("J71").Value = 4999
$ws.Range("L71").Value = 44991
$ws.Range("N71").Value = -53103

# Sheet CUL Row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 304.81818
$ws.Range("J92").Value = 341.85715
$ws.Range("L92").Value = 1025.57145
$ws.Range("N92").Value = -3521.57145

# Sheet CUL Row 111
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 7494
$ws.Range("J111").Value = 10999
$ws.Range("L111").Value = 32997
$ws.Range("N111").Value = -39131

# Sheet CUL Row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 11261.333
$ws.Range("I126").Value = 11261.333
$ws.Range("K126").Value = 33783.999
$ws.Range("M126").Value = -28843.999

# Sheet CUL Row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5159.3
$ws.Range("I137").Value = 4498.75
$ws.Range("J137").Value = 5599.6665
$ws.Range("K137").Value = 13496.25
$ws.Range("L137").Value = 16798.9995
$ws.Range("M137").Value = -8396.25
$ws.Range("N137").Value = -26998.9995

# Sheet GSM Row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6190.6665
$ws.Range("I70").Value = 5478.8
$ws.Range("K70").Value = 5478.8
$ws.Range("M70").Value = -5208.8

# Sheet GSM Row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6190.6665
$ws.Range("I73").Value = 5478.8
$ws.Range("K73").Value = 5478.8
$ws.Range("M73").Value = -4542.8

# Sheet GSM Row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3321.4
$ws.Range("I102").Value = 3209.8333
$ws.Range("K102").Value = 3209.8333
$ws.Range("M102").Value = -1587.8333

# Sheet LTW Row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4423.8237
$ws.Range("I40").Value = 4195.8667
$ws.Range("K40").Value = 4195.8667
$ws.Range("M40").Value = -4059.8667

# Sheet LTW Row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2819.5557
$ws.Range("I93").Value = 1568.8182
$ws.Range("K93").Value = 1568.8182
$ws.Range("M93").Value = -320.8181999999999

# Sheet LTW Row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2890.4
$ws.Range("I100").Value = 2350.889
$ws.Range("K100").Value = 2350.889
$ws.Range("M100").Value = -1809.889

# Sheet LTW Row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4009.35
$ws.Range("I136").Value = 1968.625
$ws.Range("J136").Value = 5369.8335
$ws.Range("K136").Value = 5905.875
$ws.Range("L136").Value = 16109.5005
$ws.Range("M136").Value = -3355.875
$ws.Range("N136").Value = -21209.5005

# Sheet WVR Row 5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1000000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1000000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1000000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1000224

# Sheet WVR Row 110
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 27644
$ws.Range("J110").Value = 27644
$ws.Range("L110").Value = 27644
$ws.Range("N110").Value = -35824

# Sheet WVR Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 25311.35
$ws.Range("I132").Value = 25783.047
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 77349.141
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -74819.141
$ws.Range("N132").Value = -21560

# Sheet WVR Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2226.25
$ws.Range("I136").Value = 1751.5625
$ws.Range("J136").Value = 4125
$ws.Range("K136").Value = 5254.6875
$ws.Range("L136").Value = 12375
$ws.Range("M136").Value = -2704.6875
$ws.Range("N136").Value = -17475
